$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Coordinador de Equipo" value (C4) ---
$ws.Range("C4").Value = "Francisco Baeza"

# --- Update activity log rows (C7:D16) ---
# Order below mirrors the order the new unique strings were introduced.
$ws.Range("C8").Value  = "Decidimos datos a pedir en los formularios"
$ws.Range("C7").Value  = "Nos reunimos via discord para ver errores de la entrega anterior"
$ws.Range("C10").Value = "Iniciamos la contruccion de los formularios"
$ws.Range("C11").Value = "Cambio barra navegadora"
$ws.Range("C12").Value = "Busqueda de la API"
$ws.Range("C13").Value = "creacion de pagina registro"
$ws.Range("C14").Value = "creacion de pagina publicar"
$ws.Range("C15").Value = "busquedad de eventos para mouse y pagina oscura"
$ws.Range("C9").Value  = "Buscar paginas donde haya referencias para la pagina "

$ws.Range("D7").Value  = "Francisco Baeza y Alexis Rodriguez"
$ws.Range("D8").Value  = "Francisco Baeza y Alexis Rodriguez"
$ws.Range("D9").Value  = "Francisco Baeza y Alexis Rodriguez"
$ws.Range("D10").Value = "Francisco Baeza "
$ws.Range("D11").Value = "Francisco Baeza y Alexis Rodriguez"
$ws.Range("D12").Value = "Alexis Rodriguez"
$ws.Range("D13").Value = "Alexis Rodriguez"
$ws.Range("D14").Value = "Francisco Baeza y Alexis Rodriguez"
$ws.Range("D15").Value = "Francisco Baeza y Alexis Rodriguez"
$ws.Range("D16").Value = "Francisco Baeza "

$ws.Range("C16").Value = "Redaccion de la bitacora"

# --- Column C width widened ---
$ws.Columns.Item(3).ColumnWidth = 57.5

# --- Active selection moved to C9 ---
$ws.Range("C9").Select()
